{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" and\n// \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" paragraphs (plus the blank\n// paragraph that separated them from the requirements list above), as in\n// the Jekyll site rebuild that dropped the scraped page-chrome text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two paragraphs that must be removed by their exact text.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    jupiterIdx = i;\n  } else if (t.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx !== -1 && copyrightIdx !== -1) {\n  // Also remove the blank paragraph immediately preceding the \"Ver no\n  // Jupiter...\" paragraph, mirroring the diff which drops that separator too.\n  const blankIdx = jupiterIdx - 1;\n  if (blankIdx >= 0 && items[blankIdx].text === \"\") {\n    items[blankIdx].delete();\n  }\n  items[jupiterIdx].delete();\n  items[copyrightIdx].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" and\n# \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" paragraphs, along with the blank\n# paragraph that separated them from the requirements list above \u2014 mirrors\n# the Jekyll rebuild that dropped this scraped page-chrome text.\n$d = $word.ActiveDocument\n\n# Locate the \"LOB1018: F\u00edsica I (Requisito fraco)\" paragraph; the three\n# paragraphs to remove immediately follow it (blank line, \"Ver no\n# Jupiter...\", \"\u00a9 2020...\").\n$reqIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -eq \"LOB1018: F\u00edsica I (Requisito fraco)`r\") {\n        $reqIndex = $i\n        break\n    }\n}\n\nif ($reqIndex -ne -1) {\n    # Delete back-to-front so earlier indices stay valid as later ones vanish.\n    # Guard on paragraph text too, so a re-run against an already-edited\n    # document (or unexpected layout) is a harmless no-op instead of an error.\n    for ($k = 3; $k -ge 1; $k--) {\n        $idx = $reqIndex + $k\n        if ($idx -gt $d.Paragraphs.Count) {\n            continue\n        }\n        $t = $d.Paragraphs.Item($idx).Range.Text\n        if ($t -eq \"`r\" -or\n            $t -eq \"Ver no Jupiter Salvar em pdf Salvar em docx`r\" -or\n            $t -like \"*Contact: luizeleno@usp.br*\") {\n            $d.Paragraphs.Item($idx).Range.Delete()\n        }\n    }\n}\n"}
